# Translate artifact slot names in column A from Chinese to English so the
# calculator can be used by both Chinese and English users.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "杯" = "goblet"
    "沙" = "sands"
    "羽" = "plume"
    "花" = "flower"
    "冠" = "circlect"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}

# Reflect the author's last active selection in the saved view state.
$ws.Range("I19").Select()
